$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Carbon Dioxide" rows (Texas block row 2 and Gulf Coast block,
# originally row 22 but which becomes row 21 after the first delete shifts rows up).
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(21).Delete()

# Insert a new "Ethylene" row at the end of the Texas block (after the existing
# Texas rows, before the Permian Basin block which now starts at row 11).
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Ethylene"
$ws.Range("B11").Value = 261.0
$ws.Range("C11").Value = 847.0
$ws.Range("D11").Value = "Texas"

# Append a new "VOC - Unclassified" row at the end of the Gulf Coast block (end of table).
$ws.Range("A31").Value = "VOC - Unclassified"
$ws.Range("B31").Value = 257.0
$ws.Range("C31").Value = 254.0
$ws.Range("D31").Value = "Gulf Coast"
